# Applies the DaCapo ZGC jython heap-8G docx edit:
#  - rows 1-3 ("100","0","58") become "0M"
#  - 10 new single-value rows are inserted right after (old) row 3
#  - the three multi-tab summary rows near the end collapse down to the
#    single values "100", "0", "58" respectively

$d = $word.ActiveDocument
$t = $d.Tables(1)

# 1) First three rows: "100" -> "0M", "0" -> "0M", "58" -> "0M"
$t.Rows(1).Cells(1).Range.Text = "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"

# 2) Insert 10 new rows right before the (old) 4th row, carrying the new values.
#    Each Rows.Add(anchor) call inserts immediately above the anchor row, so
#    adding repeatedly against the same anchor stacks rows in reverse order;
#    walk the desired values backwards to land them in forward order.
$newValues = @("14", "0.00002", "0.00008", "0.00004", "0.00001", "0.00004", "0.00005", "0.00006", "0.00065", "100.0")

$anchorRow = $t.Rows(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells(1).Range.Text = $newValues[$i]
}

# 3) Collapse the three trailing multi-tab rows into single values.
#    After inserting 10 rows above, the old row 34/35/36 are now 44/45/46.
$t.Rows(44).Cells(1).Range.Text = "100"
$t.Rows(45).Cells(1).Range.Text = "0"
$t.Rows(46).Cells(1).Range.Text = "58"
